$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean: wipe the previously used range (old table was A1:J2) so
# stale columns (G:J) and their formatting disappear entirely.
$ws.Range("A1:J6").Clear()

# ---- Header row (row 1): text labels in columns A-F ----
$ws.Range("A1").Value = "L"
$ws.Range("B1").Value = "GUA"

# Force "32" to be stored as text (matches the shared-string <t>32</t> in
# the target, not a numeric 32) - pre-format the cell as Text before typing.
$ws.Range("C1").NumberFormat = "@"
$ws.Range("C1").Value = "32"

$ws.Range("D1").Value = "Water"
$ws.Range("E1").Value = "H2"
$ws.Range("F1").Value = "CO"

# Re-apply the bold / thin-border / center-top header formatting (style
# index 1 in the original workbook) to the new A1:F1 header cells.
$header = $ws.Range("A1:F1")
$header.VerticalAlignment = -4160
$header.HorizontalAlignment = -4108
$header.Font.Bold = $true
$header.Borders.LineStyle = 1

# Column A data cells also carried that same style in the original sheet.
$colA = $ws.Range("A2:A6")
$colA.VerticalAlignment = -4160
$colA.HorizontalAlignment = -4108
$colA.Font.Bold = $true
$colA.Borders.LineStyle = 1

# ---- Data rows 2-6, columns A-F ----
$data = @(
    @(0, 0, 0, 2353.529330162142,  272561.5825838298,  391835.4054309646),
    @(0, 0, 0, 2513.069710310881,  100908.4278142457,  216842.3500618415),
    @(0, 0, 0, 2513.069710310881,  335747.7936134483,  204189.8259181912),
    @(0, 0, 0, 2513.069710310881,  44716.4143642417,   198496.9576079197),
    @(0, 0, 0, 2513.069710310881,  44295.81448838354,  168471.5104367869)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
